# Auto commit at 2026-02-14  9:28:27.23
#
# Updates the monthly "Metrics" figures (rows 2-13, column B) with the
# latest numbers. The "today" sheet pulls these values in via formulas
# (B11:B22 = Metrics!B2:B13, with E/F columns derived from those), so
# updating Metrics lets the dependent formulas on "today" recalculate
# automatically. Finally the two sheets' last-used selection (cursor
# position) is restored to match where the editor left off.

$wb = $excel.ActiveWorkbook
$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsToday   = $wb.Worksheets.Item("today")

$wsMetrics.Range("B2").Value  = 204338.46
$wsMetrics.Range("B3").Value  = 185732.49
$wsMetrics.Range("B4").Value  = 70372.91
$wsMetrics.Range("B5").Value  = 8167
$wsMetrics.Range("B6").Value  = 784974.23
$wsMetrics.Range("B7").Value  = 638526.93999999994
$wsMetrics.Range("B8").Value  = 234807.47
$wsMetrics.Range("B9").Value  = 31643
$wsMetrics.Range("B10").Value = 34886225.950000003
$wsMetrics.Range("B11").Value = 32684519.729999997
$wsMetrics.Range("B12").Value = 12180621.33
$wsMetrics.Range("B13").Value = 1349550

# Restore the last-selected cell on each sheet. Selecting a range makes
# that sheet the active one, so set Metrics first and finish on "today"
# so it stays the active tab (as it was originally).
$wsMetrics.Range("E28").Select()
$wsToday.Range("E9").Select()
